# Auto-generated edit script: updates Price (D) and Volume(1h) (E) columns
# for the cryptos list, matching the upstream commit's scraped refresh.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D holds plain-text numbers (e.g. "0.9998"). Assigning a numeric-looking
# string straight to .Value lets Excel auto-sniff it into a real Number, which
# would change the stored cell type. Flip the NumberFormat to Text ("@") for the
# affected contiguous blocks first so the assignment is kept as text, then strip
# the NumberFormat back off (ClearFormats) so the cell keeps its original (default)
# style - only its text content changes, exactly like the source diff.
$ws.Range("D4:D14").NumberFormat = "@"
$ws.Range("D16:D22").NumberFormat = "@"
$ws.Range("D25:D27").NumberFormat = "@"
$ws.Range("D29:D51").NumberFormat = "@"

$ws.Range("D2").Value = '26.887.42'
$ws.Range("E2").Value = '  -3.85%  '
$ws.Range("D3").Value = '1.729.74'
$ws.Range("E3").Value = '  -2.21%  '
$ws.Range("D4").Value = '0.9998'
$ws.Range("D5").Value = '309.38'
$ws.Range("E5").Value = '  -5.97%  '
$ws.Range("D6").Value = '0.9985'
$ws.Range("E6").Value = '  -0.16%  '
$ws.Range("D7").Value = '0.4912'
$ws.Range("E7").Value = '  +5.10%  '
$ws.Range("D8").Value = '0.3509'
$ws.Range("E8").Value = '  -0.28%  '
$ws.Range("D9").Value = '43.01'
$ws.Range("E9").Value = '  -2.03%  '
$ws.Range("D10").Value = '0.07247'
$ws.Range("E10").Value = '  -1.97%  '
$ws.Range("D11").Value = '1.052'
$ws.Range("E11").Value = '  -3.02%  '
$ws.Range("D12").Value = '0.9975'
$ws.Range("E12").Value = '  -0.32%  '
$ws.Range("D13").Value = '19.89'
$ws.Range("E13").Value = '  -3.73%  '
$ws.Range("D14").Value = '5.873'
$ws.Range("D15").Value = '1.722.67'
$ws.Range("E15").Value = '  -2.44%  '
$ws.Range("D16").Value = '6.811'
$ws.Range("E16").Value = '  -5.43%  '
$ws.Range("D17").Value = '86.74'
$ws.Range("E17").Value = '  -6.05%  '
$ws.Range("D18").Value = '0.00001035'
$ws.Range("E18").Value = '  -1.98%  '
$ws.Range("D19").Value = '0.06399'
$ws.Range("E19").Value = '  -0.32%  '
$ws.Range("D20").Value = '0.9992'
$ws.Range("E20").Value = '  -0.07%  '
$ws.Range("D21").Value = '16.53'
$ws.Range("E21").Value = '  -2.33%  '
$ws.Range("D22").Value = '5.721'
$ws.Range("E22").Value = '  -1.40%  '
$ws.Range("D23").Value = '26.977.40'
$ws.Range("E23").Value = '  -3.71%  '
$ws.Range("E24").Value = '  -1.79%  '
$ws.Range("D25").Value = '2.053'
$ws.Range("E25").Value = '  -4.78%  '
$ws.Range("D26").Value = '154.26'
$ws.Range("E26").Value = '  -5.88%  '
$ws.Range("D27").Value = '19.93'
$ws.Range("E27").Value = '  -0.68%  '
$ws.Range("D28").Value = '1.916.02'
$ws.Range("E28").Value = '  -2.69%  '
$ws.Range("D29").Value = '2.076'
$ws.Range("E29").Value = '  -5.81%  '
$ws.Range("D30").Value = '120.18'
$ws.Range("E30").Value = '  -2.62%  '
$ws.Range("D31").Value = '1.049'
$ws.Range("E31").Value = '  -2.61%  '
$ws.Range("D32").Value = '0.09287'
$ws.Range("E32").Value = '  -0.45%  '
$ws.Range("D33").Value = '3.577'
$ws.Range("E33").Value = '  -2.18%  '
$ws.Range("D34").Value = '5.379'
$ws.Range("E34").Value = '  -3.27%  '
$ws.Range("D35").Value = '0.05897'
$ws.Range("E35").Value = '  -3.58%  '
$ws.Range("D36").Value = '0.02181'
$ws.Range("E36").Value = '  -3.87%  '
$ws.Range("D37").Value = '1.429'
$ws.Range("E37").Value = '  -1.12%  '
$ws.Range("D38").Value = '10.98'
$ws.Range("E38").Value = '  -6.08%  '
$ws.Range("D39").Value = '4.750'
$ws.Range("E39").Value = '  -3.33%  '
$ws.Range("D40").Value = '0.1987'
$ws.Range("E40").Value = '  -4.13%  '
$ws.Range("D41").Value = '0.9983'
$ws.Range("E41").Value = '  -0.16%  '
$ws.Range("D42").Value = '0.5984'
$ws.Range("E42").Value = '  -2.93%  '
$ws.Range("D43").Value = '1.112'
$ws.Range("E43").Value = '  -6.92%  '
$ws.Range("D44").Value = '7.425'
$ws.Range("E44").Value = '  -4.61%  '
$ws.Range("D45").Value = '12.80'
$ws.Range("E45").Value = '  -2.56%  '
$ws.Range("D46").Value = '3.573'
$ws.Range("E46").Value = '  -4.67%  '
$ws.Range("D47").Value = '0.5609'
$ws.Range("E47").Value = '  -3.44%  '
$ws.Range("D48").Value = '119.41'
$ws.Range("E48").Value = '  -3.75%  '
$ws.Range("D49").Value = '1.838'
$ws.Range("E49").Value = '  -5.06%  '
$ws.Range("D50").Value = '0.06651'
$ws.Range("E50").Value = '  -2.39%  '
$ws.Range("D51").Value = '1.095'
$ws.Range("E51").Value = '  -2.92%  '

# Drop the temporary Text number-format again so cells end up with no explicit
# style override (matches the original workbook, which left these cells unstyled).
$ws.Range("D4:D14").ClearFormats()
$ws.Range("D16:D22").ClearFormats()
$ws.Range("D25:D27").ClearFormats()
$ws.Range("D29:D51").ClearFormats()
